# Update "想去人数" (want-to-go count) values in the 苏州-漫展信息 workbook.
# Sheet "展览" (Exhibition)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 769
$ws1.Range("F6").Value = 123
$ws1.Range("F8").Value = 126
$ws1.Range("F9").Value = 325
$ws1.Range("F10").Value = 434
$ws1.Range("F12").Value = 134
$ws1.Range("F13").Value = 11445
$ws1.Range("F14").Value = 5377

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 99

# Sheet "全部类型" (All Types) - aggregated view of the other sheets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 769
$ws4.Range("F5").Value = 99
$ws4.Range("F8").Value = 123
$ws4.Range("F10").Value = 126
$ws4.Range("F11").Value = 325
$ws4.Range("F12").Value = 434
$ws4.Range("F14").Value = 134
$ws4.Range("F15").Value = 11445
$ws4.Range("F17").Value = 5377
